# Refresh the main GSC export data: the breadcrumbs export advanced by one
# day, so the oldest date row (2025-10-14) drops off the front of the
# "Chart" table and every remaining row shifts up by one. Deleting the
# worksheet row (rather than rewriting each cell) naturally re-flows the
# dates/values and shrinks the used range from A1:C91 to A1:C90.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Row 2 holds the oldest day (2025-10-14). Deleting it shifts rows 3..91
# up into rows 2..90, so e.g. old C3 (49) becomes new C2, etc., and the
# table ends up one row shorter (A1:C90).
$ws.Rows.Item(2).Delete()
